$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I/J numeric data for rows 2 through 26
$data = @(
    @(1,4),
    @(4,6),
    @(6,8),
    @(6,8),
    @(5,7),
    @(6,7),
    @(4,7),
    @(1,5),
    @(4,6),
    @(1,5),
    @(1,6),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,6),
    @(2,7),
    @(6,9),
    @(6,9),
    @(6,8),
    @(6,9),
    @(4,7),
    @(4,6),
    @(1,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
